$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.599.09"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.563.23"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'210.49"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'0.488"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "'25.10"
$ws.Range("E8").Value = "  +5.70%  "
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'0.0587"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.786.20"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "1.565.94"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").Value = "28.606.01"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'0.515"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "'3.64"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "'61.29"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "'229.42"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "'0.998"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'3.90"
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").Value = "'8.99"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").Value = "'151.00"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("D29").Value = "'6.23"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'0.0461"
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("E31").Value = "  -2.92%  "
$ws.Range("D32").Value = "'3.17"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "1.387.08"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "'2.98"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  -3.14%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").Value = "'2.70"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").Value = "'2.29"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "'0.0162"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").Value = "'0.519"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'0.773"
$ws.Range("E43").Value = "  -1.75%  "
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "'63.91"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("E46").Value = "  -2.46%  "
$ws.Range("D47").Value = "1.699.28"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("D48").Value = "'0.868"
$ws.Range("E48").Value = "  -5.49%  "
$ws.Range("D49").Value = "'85.17"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'43.19"
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("E51").Value = "  -0.49%  "
